# Generate Report for Handback
# The handback-report generator marked the second tracked file
# (39613f77-56ff-4866-a0b3-591de88e5561.md) as handed back / in sync,
# refreshed its "Latest Handback DateTime" timestamps, and cleared the
# stale "handback file is not the latest" error for both the zh-cn and
# de-de localization targets. The Overview sheet mirrors the new status.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$status = "Handed back: in sync with en-US"

# Overview sheet: row 3 holds the 39613f77-... file; both the zh-cn and
# de-de status columns move from "Ready for handoff" to the handed-back
# status (the "Latest HO Xliff Generate Date" column is unchanged).
$overview.Range("E3").Value = $status
$overview.Range("F3").Value = $status

# zh-cn sheet: row 3 is the same file. Status, Latest Handback DateTime,
# and Error Detail all refresh.
$zhcn.Range("C3").Value = $status
$zhcn.Range("K3").Value = "2016-08-15 18:43:54"
$zhcn.Range("P3").Value = ""

# de-de sheet: same change, with its own handback timestamp.
$dede.Range("C3").Value = $status
$dede.Range("K3").Value = "2016-08-15 18:44:03"
$dede.Range("P3").Value = ""

# The Error Detail column (P) no longer holds the long "not latest"
# message in either language sheet, so Excel auto-fits it narrower.
$zhcn.Columns.Item(16).ColumnWidth = 13.7470528738839
$dede.Columns.Item(16).ColumnWidth = 13.7470528738839
